$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(0.7117080123049675, 0.6804301686572956, 0.6614864130263527, 0.653832932145491, 0.6525660976209906, 0.6613829264105391, 0.7008696152832385, 0.7803502673318405, 0.8399681022684717, 0.8673506231752981, 0.8777568094152741, 0.8755140151279761, 0.8682060082526561, 0.863734444385301, 0.8381838067141416, 0.8225759663021392, 0.8136234561363267, 0.8105965537372128, 0.8242348953207568, 0.8703515476072141, 0.9007069181457155, 0.884486189294222, 0.8234848291547507, 0.7586320231602315)
$Cvals = @(0.1600642082822716, 0.1594244198743837, 0.1590338788330001, 0.1588753248073473, 0.1588490334886679, 0.1590317380851118, 0.1598431443660928, 0.1614517753221634, 0.1626434518755957, 0.1631875420352245, 0.1633938447708942, 0.1633494021790085, 0.1632045094283896, 0.1631157927340894, 0.1626079327618939, 0.1622968736160217, 0.1621181490839589, 0.1620576689656446, 0.1623299670249665, 0.1632470608256824, 0.1638479859151758, 0.1635271250495052, 0.1623150051594351, 0.1610148149816837)
$Dvals = @(0.104985907098623, 0.1037915478909213, 0.1030930743047591, 0.1028172593192167, 0.1027719945761305, 0.1030893187929607, 0.1045668851633934, 0.1077390394663169, 0.1102346882391743, 0.1114053863090518, 0.1118537450047725, 0.1117569595739312, 0.1114421723211194, 0.1112500109556009, 0.1101588885681934, 0.1094985539800817, 0.1091220835030882, 0.1089951917011831, 0.1095685027830271, 0.1115344965759419, 0.1128487445072324, 0.1121446365596483, 0.1095368690681155, 0.1068517296979508)
$Evals = @(0.1422608133583481, 0.1427241575655867, 0.1430468559831795, 0.1431879886222873, 0.1432120059459656, 0.1430487203220494, 0.1424126596264745, 0.141467453508529, 0.1409558612950743, 0.1407625532655263, 0.1406949975765031, 0.1407092961569436, 0.1407568823945695, 0.1407867649351218, 0.1409692857030578, 0.1410913361946253, 0.1411652479356054, 0.1411909113221164, 0.141077959744921, 0.1407427521138285, 0.1405565726509863, 0.1406529373274967, 0.1410839955764107, 0.1416909511064937)
$Fvals = @(2.307636691385504, 2.317594672145383, 2.324613667606336, 2.327701745878969, 2.328228284332887, 2.32465439187424, 2.310882582939946, 2.291043875853738, 2.280824541167533, 2.277118650276392, 2.275850674554832, 2.276117739063423, 2.277011621546805, 2.277576772090342, 2.281085684611, 2.28347963475278, 2.284945338755406, 2.285456853760891, 2.283215609253929, 2.276745394651627, 2.273305638846779, 2.27506939190468, 2.283334696718711, 2.295644810000852)
$Jvals = @(0.213893841753503, 0.2151499831879118, 0.2159684523112748, 0.2163138712992883, 0.2163719463156628, 0.2159730626003782, 0.2143171758988007, 0.2114435662719991, 0.2095588293746324, 0.2087503334748426, 0.2084511867429715, 0.2085153016844998, 0.2087255820300093, 0.2088552976617084, 0.2096126487255798, 0.2100897672085686, 0.2103687949289235, 0.2104640597668528, 0.2100385009788237, 0.2086636273865938, 0.2078059385938023, 0.2082599687742235, 0.2100616637191042, 0.212181088280289)
$Kvals = @(0.3364172470825508, 0.3073175970681348, 0.2895201197247843, 0.2822855278896697, 0.2810853318375024, 0.2894224778879959, 0.3263694656212976, 0.3993596803011599, 0.4532961433037315, 0.4778974195050978, 0.487222278572176, 0.4852136159301779, 0.4786644059783214, 0.4746539686144047, 0.4516896633127772, 0.4376181755847028, 0.4295308050772917, 0.4267936366065044, 0.4391154738117109, 0.4805878304242697, 0.5077439969854822, 0.4932456852503719, 0.438438537411912, 0.3795581356595221)
$Lvals = @(0.2358674811642132, 0.2341666889500544, 0.2332019661925457, 0.232828928303249, 0.2327682018290389, 0.2331968537969473, 0.2352645836107996, 0.2399473750425329, 0.2437666945461245, 0.2455856696069034, 0.2462861162180872, 0.2461347461754642, 0.2456430629764412, 0.2453434061675921, 0.2436494538186835, 0.2426310937849223, 0.2420530408014017, 0.2418586438029422, 0.2427387056431485, 0.2457871670240053, 0.247847321097332, 0.2467416020525448, 0.2426900312167106, 0.2386137366509828)
$Mvals = @(0.2029199323890616, 0.1974650344182045, 0.1941925677695941, 0.1928784731158828, 0.1926614481192424, 0.1941747664497235, 0.2010232090928525, 0.2150578027648429, 0.2257321894802544, 0.2306660087951968, 0.2325454139542273, 0.2321401600886119, 0.230820407430528, 0.2300134589564138, 0.2254113101461783, 0.2226079158844172, 0.2210028265078421, 0.2204606372563482, 0.2229055824229604, 0.231207751319964, 0.2366981868352624, 0.2337619834413402, 0.2227709867522236, 0.2111968977022158)
$Ovals = @(6.287721418772492, 6.321571816392662, 6.344846736307005, 6.354958314533206, 6.356675211878525, 6.344980565708937, 6.298876441221751, 6.228205355032458, 6.188288291082785, 6.172729985771497, 6.167211865722123, 6.168383688126227, 6.172268524218907, 6.17469672294672, 6.189357348294891, 6.199016827186, 6.204817469531747, 6.206823519164743, 6.197963229949494, 6.171117321160949, 6.155748653429157, 6.163752175038496, 6.198438790961262, 6.245213508256796)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $Bvals[$i]
    $ws.Cells.Item($row, 3).Value = $Cvals[$i]
    $ws.Cells.Item($row, 4).Value = $Dvals[$i]
    $ws.Cells.Item($row, 5).Value = $Evals[$i]
    $ws.Cells.Item($row, 6).Value = $Fvals[$i]
    $ws.Cells.Item($row, 10).Value = $Jvals[$i]
    $ws.Cells.Item($row, 11).Value = $Kvals[$i]
    $ws.Cells.Item($row, 12).Value = $Lvals[$i]
    $ws.Cells.Item($row, 13).Value = $Mvals[$i]
    $ws.Cells.Item($row, 15).Value = $Ovals[$i]
}

Write-Output "done"